# Update "想去人数" (people interested) figures on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 9176
$ws1.Range("F4").Value = 475
$ws1.Range("F5").Value = 456

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 9176
$ws4.Range("F4").Value = 475
$ws4.Range("F6").Value = 456
